# Generate Report for Archive
# - Flip the localization status shown in the report from
#   "Ready for handoff" to "In Translation" on every sheet that shows it.
# - Re-narrow the Status/locale columns now that the new text is shorter
#   than "Ready for handoff" (mirrors an auto-fit-to-content pass).

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: zh-cn / de-de status cells (columns E & F, row 2) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# Columns E & F were sized to fit "Ready for handoff"; shrink them now that
# the status text is shorter.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

# --- Per-locale detail sheets: Status column (column C, row 2) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Columns.Item(3).ColumnWidth = 12.5
